# "Full-screen slides on PAM_STIM_2024"
# Renumber the Diapositive (slide) references in the "Img" column so the
# story's slide deck lines up with the new full-screen slide numbering:
#   Diapositive6b  -> Diapositive7
#   Diapositive6c  -> Diapositive8
#   Diapositive10  -> Diapositive12
#   Diapositive12  -> Diapositive14
#   Diapositive13  -> Diapositive15
#   Diapositive11  -> Diapositive13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = "Diapositive7"
$ws.Range("C9").Value  = "Diapositive8"
$ws.Range("C15").Value = "Diapositive12"
$ws.Range("C16").Value = "Diapositive14"
$ws.Range("C17").Value = "Diapositive15"
$ws.Range("C18").Value = "Diapositive13"

# Leave the cursor where the author left it when saving.
$null = $ws.Range("C19").Select()
